$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# HP column type changed from int to float
$ws.Range("D2").Value = "float"

# MP (column C) starting value changed from 100 to 0 for all hero rows
for ($r = 11; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
}

# Rename "Earch*" typo to "Earth*"
$ws.Range("A15").Value = "EarthSaber"
$ws.Range("H15").Value = "EarthSaber"
$ws.Range("A20").Value = "EarthArcher"
$ws.Range("H20").Value = "EarthArcher"
$ws.Range("A25").Value = "EarthWizard"
$ws.Range("H25").Value = "EarthWizard"

# Update selection / view state
$null = $ws.Range("D27").Select()
